$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K header (2022) - copy formatting from J4 then set value
$ws.Range("J4").Copy($ws.Range("K4"))
$ws.Range("K4").Value = 2022

# Row 5 - copy formatting from J5 then set value
$ws.Range("J5").Copy($ws.Range("K5"))
$ws.Range("K5").Value = 26.495524312074597

# Row 6
$ws.Range("J6").Copy($ws.Range("K6"))
$ws.Range("K6").Value = 59.383769502755833

# Row 7
$ws.Range("J7").Copy($ws.Range("K7"))
$ws.Range("K7").Value = 38.32334404557426

# Row 8
$ws.Range("J8").Copy($ws.Range("K8"))
$ws.Range("K8").Value = 48.136790950525594

# Row 9
$ws.Range("J9").Copy($ws.Range("K9"))
$ws.Range("K9").Value = 46.63213064070051

# Row 10
$ws.Range("J10").Copy($ws.Range("K10"))
$ws.Range("K10").Value = 32.657429481680126

# Row 11
$ws.Range("J11").Copy($ws.Range("K11"))
$ws.Range("K11").Value = 31.457245964894081

# Row 12
$ws.Range("J12").Copy($ws.Range("K12"))
$ws.Range("K12").Value = 22.734405597714229

# Row 13
$ws.Range("J13").Copy($ws.Range("K13"))
$ws.Range("K13").Value = -0.19691879995369213

# Row 14
$ws.Range("J14").Copy($ws.Range("K14"))
$ws.Range("K14").Value = 33.158040409631916

# Update the active-cell selection shown in the saved sheet view
$ws.Range("M7").Select()
